# Fix the "Speed Up" column formula in Table1 (Sheet1!E3:E9).
# Previously the formula computed D<row>/$D$3 (slow-down relative to the
# single-thread baseline); it should instead compute $D$3/D<row> (true
# speed-up factor, i.e. baseline time divided by the row's time).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

for ($r = 3; $r -le 9; $r++) {
    $ws.Cells.Item($r, 5).Formula = "=`$D`$3/D$r"
}

# Recalculate so the cached <v> values stored in the file reflect the
# corrected formula.
$wb.Application.Calculate()

# Restore the last selected cell (the author had clicked E20 before saving).
$ws.Range("E20").Select() | Out-Null
